# Se ajusta funcionalidad para mostrar las informacion segun los filtros
#
# The "departamentos" list contained a duplicated/typo entry:
#   row 14 -> "Valle de Cauca "   (typo)
#   row 15 -> "Valle del Cauca "  (correct)
# Remove the erroneous row (row 14) so the correct one shifts up and
# becomes the new (and last) row 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate/incorrect "Valle de Cauca " row, shifting the
# remaining rows (including the correct "Valle del Cauca ") up.
$ws.Rows(14).Delete()

# Leave the selection on the last data row, matching the saved view state.
$ws.Range("A14").Select()
